$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'69.644.66"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = "'2.493.13"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = "'568.54"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.00%  '
$ws.Range('D6').Value = "'165.63"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.91%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -1.38%  '
$ws.Range('D9').Value = "'2.490.60"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('E10').Value = '  -3.22%  '
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('D12').Value = "'0.355"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = "'4.92"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = "'2.948.88"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.74%  '
$ws.Range('D15').Value = "'69.527.22"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('E16').Value = '  -1.20%  '
$ws.Range('D17').Value = "'24.33"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.10%  '
$ws.Range('D18').Value = "'2.494.69"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.05%  '
$ws.Range('D19').Value = "'11.19"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').Value = "'7.38"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.06%  '
$ws.Range('D21').Value = "'346.34"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('E22').Value = '  -1.90%  '
$ws.Range('E23').Value = '  -5.45%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = "'70.38"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').Value = "'3.89"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.69%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = "'2.615.29"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').Value = "'8.66"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.29%  '
$ws.Range('D29').Value = "'0.993"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = "'7.82"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = "'0.0₃0879"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.37%  '
$ws.Range('D32').Value = "'445.76"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.41%  '
$ws.Range('D33').Value = "'1.19"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.75%  '
$ws.Range('D34').Value = "'1.00"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.27%  '
$ws.Range('E35').Value = '  -3.20%  '
$ws.Range('D36').Value = "'155.79"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.22%  '
$ws.Range('E37').Value = '  -4.53%  '
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('D39').Value = "'18.19"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('E41').Value = '  -2.27%  '
$ws.Range('D42').Value = "'4.61"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.69%  '
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('D44').Value = "'38.06"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('E45').Value = '  -7.68%  '
$ws.Range('E46').Value = '  -8.76%  '
$ws.Range('D47').Value = "'139.85"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.51%  '
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range('E49').Value = '  -4.34%  '
$ws.Range('D50').Value = "'0.0728"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('D51').Value = "'0.574"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.82%  '
